# Refactor security vulnerability checks
# Appends one new telemetry row to each of the four worksheets, matching
# new rows captured by the logging tool after the workbook was last saved.

function Set-LogRow {
    param(
        $Sheet,
        $Row,
        $TimeValue,
        $TotalLenHex,
        $IdHex,
        $ActualLenHex,
        $ChecksumHex,
        $TotalLenDec,
        $IdDec,
        $ActualLenDec,
        $ChecksumDec
    )

    # Column A: timestamp, reuse the same date/time style as the rest of the column.
    $Sheet.Cells.Item($Row, 1).Value = $TimeValue
    $Sheet.Cells.Item($Row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # Columns B-E: raw hex byte strings (stored as plain text).
    $Sheet.Cells.Item($Row, 2).Value = $TotalLenHex
    $Sheet.Cells.Item($Row, 3).Value = $IdHex
    $Sheet.Cells.Item($Row, 4).Value = $ActualLenHex
    $Sheet.Cells.Item($Row, 5).Value = $ChecksumHex

    # Column F: decoded total length.
    $Sheet.Cells.Item($Row, 6).Value = $TotalLenDec

    # Column G: decoded ID - some values overflow double precision in the
    # source log and must stay as exact-digit text instead of being
    # rounded to a float.
    if ($IdDec -is [string]) {
        $Sheet.Cells.Item($Row, 7).Value = "'" + $IdDec
        $Sheet.Cells.Item($Row, 7).Style = "Normal"
    } else {
        $Sheet.Cells.Item($Row, 7).Value = $IdDec
    }

    # Column H: decoded actual length.
    $Sheet.Cells.Item($Row, 8).Value = $ActualLenDec

    # Column I: decoded checksum.
    $Sheet.Cells.Item($Row, 9).Value = $ChecksumDec
}

$wb = $excel.ActiveWorkbook

# --- Sheet 1: ROW50-FE-LIFTER -> new row 79 -------------------------------
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$idDec1 = [double]"5.68631262647114e+23"
Set-LogRow $ws1 79 45763.7654925 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x4a" "0xe" 400 $idDec1 330 14

# --- Sheet 2: ROW50-MID-LIFTER -> new row 81 ------------------------------
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$idDec2 = "568631262647113771663628"
Set-LogRow $ws2 81 45763.72850694445 "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x4e" "0x19" 400 $idDec2 334 25

# --- Sheet 3: ROW11-FE-LIFTER -> new row 79 -------------------------------
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$idDec3 = [double]"5.68631262647114e+23"
Set-LogRow $ws3 79 45763.80038960648 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x4a" "0x14" 400 $idDec3 330 20

# --- Sheet 4: ROW11-MID-LIFTER -> new row 79 ------------------------------
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$idDec4 = [double]"5.68631262647114e+23"
Set-LogRow $ws4 79 45763.92296003472 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x52" "0x19" 400 $idDec4 338 25
